# feature: 寻址方式 switcher 改为 select
# Adds an "Addressing mode(Optional)" column (O) with Static/Dynamic option
# strings, and re-splits the rich-text run of the "Speed limit" header so the
# zero-width-space glyphs carry their own (MS Gothic) font run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-run the "Speed ​​limit M/s(Optional)" header (N1) rich text so the
#     zero-width-space characters (positions 7-8) get their own font run. ---
$chars = $ws.Range("N1").Characters(7, 2)
$chars.Font.Name = "MS Gothic"

# --- New column O: "Addressing mode(Optional)" switcher values ---
$ws.Columns.Item(15).ColumnWidth = 26.41

$ws.Range("O1").Value = "Addressing mode(Optional)"
$ws.Range("O3").Value = "Static"
$ws.Range("O2").Value = "Dynamic"

$ws.Range("G17").Select()
